# Append the new melee-weapon / creature-attack rows (46-51) that describe
# the "Dancer in Darkness" (Touch/Bite/Enthrall) and "Krashtkid" (Fangs/
# Digging Claw/Tongue) entries, matching the upstream data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# Row 46
$ws.Range("A46").Value = "Touch"
$ws.Range("B46").Value = 90
$ws.Range("E46").Value = 9
$ws.Range("F46").Value = "Dancer in Darkness"
$ws.Range("G46").Value = "Touch"
$ws.Range("H46").Value = "Touch: A victim hit by a Dancer in Darkness’ touch must match their magic points against that of the vampire, whether the touch penetrated armor or not. If the vampire wins the contest, the victim gives 1D4 magic points to the vampire."
$ws.Range("J46").Value = "Special"

# Row 47
$ws.Range("A47").Value = "Bite"
$ws.Range("B47").Value = 90
$ws.Range("C47").Value = "1D3"
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = "Dancer in Darkness"
$ws.Range("G47").Value = "Bite"
$ws.Range("H47").Value = "Bite: The bite of a Dancer in Darkness does damage equal to its damage bonus, or 1D3 in any case. If the bite penetrates armor, the vampire will stay attached and take 1D6 hit points (blood) from the victim—increasing the vampire’s hit points in the process."
$ws.Range("J47").Value = "+special"

# Row 48
$ws.Range("A48").Value = "Enthrall"
$ws.Range("B48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = "Dancer in Darkness"
$ws.Range("J48").Value = "Auto CHA vs. POW"

# Row 49
$ws.Range("A49").Value = "Fangs"
$ws.Range("B49").Value = 50
$ws.Range("C49").Value = "2D6"
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = "Krashtkid"
$ws.Range("G49").Value = "Fangs"
$ws.Range("H49").Value = "Fangs: When a bite does damage, the poison’s POT (equal to the krarshtkid’s CON) must try to overcome the target’s CON. If successful, the target will be paralyzed for 20–CON days. Mineral antidote is half-effective against this venom."

# Row 50
$ws.Range("A50").Value = "Digging Claw"
$ws.Range("B50").Value = 50
$ws.Range("C50").Value = "1D6"
$ws.Range("E50").Value = 6
$ws.Range("F50").Value = "Krashtkid"

# Row 51
$ws.Range("A51").Value = "Tongue"
$ws.Range("B51").Value = 75
$ws.Range("E51").Value = 6
$ws.Range("F51").Value = "Krashtkid"
$ws.Range("G51").Value = "Tongue"
$ws.Range("H51").Value = "Tongue: A hit by the tongue attacks the POW of the target with a resistance roll. If successful, the target will take 4D3 damage as if from a Disruption spell."

# Columns F (Creature), H (SpecialText) and J (Notes) now hold longer text
# (creature names and the new multi-sentence special-attack descriptions),
# so widen them the same way the source workbook does.
$ws.Range("F1").ColumnWidth = 19.620946884155273
$ws.Range("H1").ColumnWidth = 238.189697265625
$ws.Range("J1").ColumnWidth = 20.392431259155273

